$d = $word.ActiveDocument

# Step 1: locate the paragraph "Overall App design and navigation" inside its
# table cell and append a new (empty) paragraph right after it, within the
# same cell.
$rng = $d.Content
$rng.Find.Execute("Overall App design and navigation", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$cell = $rng.Cells(1)
$cellRange = $cell.Range
$cellRange.InsertParagraphAfter()

# Step 2: re-locate the same cell (fresh Find/Range so the newly inserted
# paragraph is visible) and put the new text into the second paragraph of
# that cell.
$rng2 = $d.Content
$rng2.Find.Execute("Overall App design and navigation", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$cell2 = $rng2.Cells(1)
$cellRange2 = $cell2.Range
$newPara = $cellRange2.Paragraphs(2)
$newPara.Range.InsertAfter("Overall code quality")
